$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")

$ws.Range("B7").Value = "Good"
$ws.Range("C7").Value = "Leading group meetings, using git to upload and good at sharing "

$ws.Range("B19").Value = "Excellent"
$ws.Range("C19").Value = "Good at communicating, active using Jens"

$ws.Range("C9").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
